# Generate Report for Handback
# Updates the handoff/handback timestamps for the ad764da2-... file that was
# just handed back, across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the ad764da2-... file; refresh its
# "Latest HO Xliff Generate Date" column (G).
$wsOverview.Range("G3").Value = "2016-08-13 07:00:38"

# zh-cn sheet: row 3 is the ad764da2-... file; refresh its
# "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K).
$wsZhCn.Range("H3").Value = "2016-08-13 07:00:29"
$wsZhCn.Range("K3").Value = "2016-08-13 07:00:57"

# de-de sheet: row 3 is the ad764da2-... file; refresh its
# "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K).
$wsDeDe.Range("H3").Value = "2016-08-13 07:00:38"
$wsDeDe.Range("K3").Value = "2016-08-13 07:01:11"
